# Auto-generated edit script to apply numeric updates from the commit diff
# to the Ixion_Profits workbook (multi-sheet: ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 76
$ws.Range("H76").Value = 28686.75
$ws.Range("I76").Value = 55771.895
$ws.Range("J76").Value = 4181.143
$ws.Range("K76").Value = 55771.895
$ws.Range("L76").Value = 4181.143
$ws.Range("M76").Value = -55456.895
$ws.Range("N76").Value = -4811.143
# Row 79
$ws.Range("H79").Value = 28686.75
$ws.Range("I79").Value = 55771.895
$ws.Range("J79").Value = 4181.143
$ws.Range("K79").Value = 55771.895
$ws.Range("L79").Value = 4181.143
$ws.Range("M79").Value = -54679.895
$ws.Range("N79").Value = -6365.143
# Row 125
$ws.Range("H125").Value = 7376.4
$ws.Range("I125").Value = 8294
$ws.Range("J125").Value = 6000
$ws.Range("K125").Value = 74646
$ws.Range("L125").Value = 54000
$ws.Range("M125").Value = -72186
$ws.Range("N125").Value = -58920

$ws = $wb.Worksheets.Item("ARM")
# Row 33
$ws.Range("H33").Value = 9000
$ws.Range("I33").Value = 9000
$ws.Range("J33").Value = 9000
$ws.Range("K33").Value = 9000
$ws.Range("L33").Value = 9000
$ws.Range("M33").Value = -8671
$ws.Range("N33").Value = -9658
# Row 61
$ws.Range("H61").Value = 4780.222
$ws.Range("I61").Value = 6581.4546
$ws.Range("J61").Value = 3541.875
$ws.Range("K61").Value = 6581.4546
$ws.Range("L61").Value = 3541.875
$ws.Range("M61").Value = -6369.4546
$ws.Range("N61").Value = -3965.875
# Row 74
$ws.Range("H74").Value = 2077.64
$ws.Range("I74").Value = 1729.2667
$ws.Range("J74").Value = 2600.2
$ws.Range("K74").Value = 1729.2667
$ws.Range("L74").Value = 2600.2
$ws.Range("M74").Value = -855.2666999999999
$ws.Range("N74").Value = -4348.2
# Row 77
$ws.Range("H77").Value = 2077.64
$ws.Range("I77").Value = 1729.2667
$ws.Range("J77").Value = 2600.2
$ws.Range("K77").Value = 8646.333499999999
$ws.Range("L77").Value = 13001
$ws.Range("M77").Value = -4278.333499999999
$ws.Range("N77").Value = -21737
# Row 132
$ws.Range("H132").Value = 2954.9575
$ws.Range("I132").Value = 1629.5172
$ws.Range("J132").Value = 5090.3887
$ws.Range("K132").Value = 4888.5516
$ws.Range("L132").Value = 15271.1661
$ws.Range("M132").Value = -2358.5516
$ws.Range("N132").Value = -20331.1661
# Row 136
$ws.Range("H136").Value = 4780.222
$ws.Range("I136").Value = 6581.4546
$ws.Range("J136").Value = 3541.875
$ws.Range("K136").Value = 19744.3638
$ws.Range("L136").Value = 10625.625
$ws.Range("M136").Value = -17194.3638
$ws.Range("N136").Value = -15725.625

$ws = $wb.Worksheets.Item("BSM")
# Row 107
$ws.Range("H107").Value = 1900
$ws.Range("I107").Value = 2200
$ws.Range("K107").Value = 2200
$ws.Range("M107").Value = -280
# Row 134
$ws.Range("H134").Value = 5991.357
$ws.Range("I134").Value = 11258.546
$ws.Range("J134").Value = 2583.1765
$ws.Range("K134").Value = 33775.638
$ws.Range("L134").Value = 7749.529500000001
$ws.Range("M134").Value = -31240.638
$ws.Range("N134").Value = -12819.5295

$ws = $wb.Worksheets.Item("CRP")
# Row 16
$ws.Range("H16").Value = 1927.3334
$ws.Range("I16").Value = 1673.6364
$ws.Range("J16").Value = 2625
$ws.Range("K16").Value = 1673.6364
$ws.Range("L16").Value = 2625
$ws.Range("M16").Value = -1386.6364
$ws.Range("N16").Value = -3199
# Row 58
$ws.Range("H58").Value = 1796.6744
$ws.Range("I58").Value = 1632.8096
$ws.Range("J58").Value = 1953.091
$ws.Range("K58").Value = 1632.8096
$ws.Range("L58").Value = 1953.091
$ws.Range("M58").Value = -1429.8096
$ws.Range("N58").Value = -2359.091
# Row 99
$ws.Range("H99").Value = 9628437
$ws.Range("I99").Value = 20285.125
$ws.Range("J99").Value = 25001480
$ws.Range("K99").Value = 20285.125
$ws.Range("L99").Value = 25001480
$ws.Range("M99").Value = -18787.125
$ws.Range("N99").Value = -25004476
# Row 107
$ws.Range("H107").Value = 844.3333
$ws.Range("I107").Value = 404.66666
$ws.Range("J107").Value = 1137.4445
$ws.Range("K107").Value = 404.66666
$ws.Range("L107").Value = 1137.4445
$ws.Range("M107").Value = 1515.33334
$ws.Range("N107").Value = -4977.4445
# Row 113
$ws.Range("H113").Value = 1927.3334
$ws.Range("I113").Value = 1673.6364
$ws.Range("J113").Value = 2625
$ws.Range("K113").Value = 1673.6364
$ws.Range("L113").Value = 2625
$ws.Range("M113").Value = 496.3635999999999
$ws.Range("N113").Value = -6965
# Row 122
$ws.Range("H122").Value = 4631067.5
$ws.Range("I122").Value = 4631067.5
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 13893202.5
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -13890752.5
$ws.Range("N122").ClearContents()
# Row 126
$ws.Range("H126").Value = 9628437
$ws.Range("I126").Value = 20285.125
$ws.Range("J126").Value = 25001480
$ws.Range("K126").Value = 60855.375
$ws.Range("L126").Value = 75004440
$ws.Range("M126").Value = -58385.375
$ws.Range("N126").Value = -75009380
# Row 136
$ws.Range("H136").Value = 1796.6744
$ws.Range("I136").Value = 1632.8096
$ws.Range("J136").Value = 1953.091
$ws.Range("K136").Value = 4898.4288
$ws.Range("L136").Value = 5859.272999999999
$ws.Range("M136").Value = -2348.4288
$ws.Range("N136").Value = -10959.273

$ws = $wb.Worksheets.Item("CUL")
# Row 21
$ws.Range("H21").Value = 1519.8
$ws.Range("I21").Value = 0
$ws.Range("J21").Value = 1519.8
$ws.Range("K21").Value = 0
$ws.Range("L21").Value = 4559.4
$ws.Range("M21").ClearContents()
$ws.Range("N21").Value = -4905.4
# Row 55
$ws.Range("H55").Value = 2755.6
$ws.Range("I55").Value = 490
$ws.Range("J55").Value = 2850
$ws.Range("K55").Value = 1470
$ws.Range("L55").Value = 8550
$ws.Range("M55").Value = -1293
$ws.Range("N55").Value = -8904
# Row 131
$ws.Range("H131").Value = 1961828.6
$ws.Range("I131").Value = 12500850
$ws.Range("J131").Value = 1080.4651
$ws.Range("K131").Value = 37502550
$ws.Range("L131").Value = 3241.3953
$ws.Range("M131").Value = -37497510
$ws.Range("N131").Value = -13321.3953

$ws = $wb.Worksheets.Item("GSM")
# Row 70
$ws.Range("H70").Value = 6039.7393
$ws.Range("I70").Value = 6389.1177
$ws.Range("J70").Value = 5049.8335
$ws.Range("K70").Value = 6389.1177
$ws.Range("L70").Value = 5049.8335
$ws.Range("M70").Value = -6119.1177
$ws.Range("N70").Value = -5589.8335
# Row 73
$ws.Range("H73").Value = 6039.7393
$ws.Range("I73").Value = 6389.1177
$ws.Range("J73").Value = 5049.8335
$ws.Range("K73").Value = 6389.1177
$ws.Range("L73").Value = 5049.8335
$ws.Range("M73").Value = -5453.1177
$ws.Range("N73").Value = -6921.8335
# Row 80
$ws.Range("H80").Value = 2402.0967
$ws.Range("I80").Value = 2350.7144
$ws.Range("K80").Value = 2350.7144
$ws.Range("M80").Value = -1352.7144
# Row 83
$ws.Range("H83").Value = 2402.0967
$ws.Range("I83").Value = 2350.7144
$ws.Range("K83").Value = 11753.572
$ws.Range("M83").Value = -6761.572
# Row 97
$ws.Range("H97").Value = 1173.5714
$ws.Range("I97").Value = 1142.5555
$ws.Range("J97").Value = 2011
$ws.Range("K97").Value = 1142.5555
$ws.Range("L97").Value = 2011
$ws.Range("M97").Value = -646.5554999999999
$ws.Range("N97").Value = -3003
# Row 113
$ws.Range("H113").Value = 83335300
$ws.Range("I113").Value = 166667860
$ws.Range("J113").Value = 2733.1667
$ws.Range("K113").Value = 166667860
$ws.Range("L113").Value = 2733.1667
$ws.Range("M113").Value = -166665690
$ws.Range("N113").Value = -7073.1667
# Row 132
$ws.Range("H132").Value = 32648.324
$ws.Range("I132").Value = 81231.08
$ws.Range("J132").Value = 2573.2856
$ws.Range("K132").Value = 243693.24
$ws.Range("L132").Value = 7719.8568
$ws.Range("M132").Value = -241163.24
$ws.Range("N132").Value = -12779.8568

$ws = $wb.Worksheets.Item("LTW")
# Row 100
$ws.Range("H100").Value = 1442.909
$ws.Range("I100").Value = 1208
$ws.Range("J100").Value = 2500
$ws.Range("K100").Value = 1208
$ws.Range("L100").Value = 2500
$ws.Range("M100").Value = -667
$ws.Range("N100").Value = -3582
# Row 136
$ws.Range("H136").Value = 5571.088
$ws.Range("I136").Value = 5081.4546
$ws.Range("J136").Value = 6244.3335
$ws.Range("K136").Value = 15244.3638
$ws.Range("L136").Value = 18733.0005
$ws.Range("M136").Value = -12694.3638
$ws.Range("N136").Value = -23833.0005

$ws = $wb.Worksheets.Item("WVR")
# Row 119
$ws.Range("H119").Value = 57725
$ws.Range("J119").Value = 57725
$ws.Range("L119").Value = 57725
$ws.Range("N119").Value = -67401
# Row 122
$ws.Range("H122").Value = 2230
$ws.Range("J122").Value = 1950
$ws.Range("L122").Value = 5850
$ws.Range("N122").Value = -10750
# Row 136
$ws.Range("H136").Value = 2747.9824
$ws.Range("I136").Value = 3456.2727
$ws.Range("J136").Value = 1774.0834
$ws.Range("K136").Value = 10368.8181
$ws.Range("L136").Value = 5322.2502
$ws.Range("M136").Value = -7818.8181
$ws.Range("N136").Value = -10422.2502

Write-Host "Applied all Sheets updates."